$d = $word.ActiveDocument

# 1. Locate the paragraph holding the old endpoint text and replace its
#    content with the new endpoint text (single run for now).
$old = "/api/clinics/{id}/doctors/**"
$new = "/api/doctors/**"

$find = $d.Content
$found = $find.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find target text '$old'"
}

# $find now spans the freshly inserted replacement text ("/api/doctors/**").
$rangeStart = $find.Start

# 2. Split the single run into four runs that share identical run
#    properties (Arial / sz 26) but are distinct <w:r> elements, matching
#    the target OOXML:  "/api/"  "doctors"  "/"  "**"
#
#    Toggling a character property on/off on a sub-range forces Word to
#    materialize a run boundary there without altering the final
#    formatting (Bold ends up back at its original value).

$docStart = $rangeStart
$offDoctors = $new.IndexOf("doctors")
$offSlash = $offDoctors + "doctors".Length
$offStars = $new.LastIndexOf("**")

# Split "doctors" off from "/api/" (run boundary before + after "doctors")
$doctorsRange = $d.Range($docStart + $offDoctors, $docStart + $offSlash)
$doctorsRange.Bold = 1
$doctorsRange.Bold = 0

# Split trailing "/" off from "**" (run boundary between them)
$slashRange = $d.Range($docStart + $offSlash, $docStart + $offStars)
$slashRange.Bold = 1
$slashRange.Bold = 0

Write-Output "done"
